$wb = $excel.ActiveWorkbook

# --- Sheet1 "babies": add a new baby bb11 housed in new room r11 ---
$wsBabies = $wb.Worksheets.Item("babies")
$wsBabies.Activate()
$wsBabies.Range("A12").Value = "bb11"
$wsBabies.Range("B12").Value = "leave_hospital"
$wsBabies.Range("C12").Value = "r11"
$wsBabies.Range("C13").Select()

# --- Sheet2 "rooms": add new room r11 + a "priority" column ---
$wsRooms = $wb.Worksheets.Item("rooms")
$wsRooms.Activate()

# Insert a new row 12 (old row 12 "out" shifts down to row 13)
$wsRooms.Rows.Item(12).Insert()

$wsRooms.Range("A12").Value = "r11"
$wsRooms.Range("C12").Value = "yes"
$wsRooms.Range("F12").Value = "neo"
$wsRooms.Range("G12").Value = 1

# New "priority" column header
$wsRooms.Range("H1").Value = "priority"

$wsRooms.Range("H1").Select()
